$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.705.61"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.890.28"
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("D5").Value = "'313.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("D7").Value = "'0.4811"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.64%  "
$ws.Range("D8").Value = "'0.3786"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.64%  "
$ws.Range("D9").Value = "'0.07331"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "'0.9186"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("D11").Value = "'20.45"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("D13").Value = "1.914.09"
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("D14").Value = "'5.470"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").Value = "'6.584"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").Value = "'90.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "'0.000008796"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").Value = "'1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").Value = "27.768.78"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("D21").Value = "'14.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("D22").Value = "'5.119"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("D23").Value = "2.170.94"
$ws.Range("E23").Value = "  +2.46%  "
$ws.Range("D25").Value = "'154.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("D26").Value = "'1.902"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.74%  "
$ws.Range("D27").Value = "'18.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").Value = "'2.105"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.33%  "
$ws.Range("D29").Value = "'116.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("D30").Value = "'4.914"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("D31").Value = "'0.08936"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").Value = "'3.144"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.57%  "
$ws.Range("D33").Value = "'1.232"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("D34").Value = "'0.7606"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("D35").Value = "'4.628"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("D36").Value = "'0.02031"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("D37").Value = "'2.525"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.17%  "
$ws.Range("D38").Value = "'1.091"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.06%  "
$ws.Range("D39").Value = "'0.05253"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.972"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5431"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.46%  "
$ws.Range("D42").Value = "'6.960"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("E43").Value = "  -0.57%  "
$ws.Range("D44").Value = "'8.297"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.50%  "
$ws.Range("D45").Value = "'109.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.59%  "
$ws.Range("D46").Value = "'10.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("D47").Value = "'0.4774"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("D48").Value = "'1.003"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("D49").Value = "'1.639"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.35%  "
$ws.Range("D50").Value = "'67.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").Value = "'0.06059"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.74%  "
